# Update "想去人数" (want-to-go count, column F) figures on the
# "展览" (Exhibitions, sheet 1) and "全部类型" (All types, sheet 4) tabs
# to reflect freshly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 4426
$ws1.Range("F3").Value  = 104
$ws1.Range("F4").Value  = 336
$ws1.Range("F7").Value  = 46
$ws1.Range("F10").Value = 313
$ws1.Range("F11").Value = 247
$ws1.Range("F12").Value = 2954
$ws1.Range("F13").Value = 147
$ws1.Range("F14").Value = 1526

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 4426
$ws4.Range("F3").Value  = 104
$ws4.Range("F4").Value  = 336
$ws4.Range("F8").Value  = 46
$ws4.Range("F11").Value = 313
$ws4.Range("F12").Value = 248
$ws4.Range("F13").Value = 2954
$ws4.Range("F14").Value = 147
$ws4.Range("F15").Value = 1526
